# Update cryptocurrency price/volume data to match the Jan 22 2023 10:56 UTC refresh.
# Values are written as text (Price/Volume(1h) columns store formatted strings like
# "302.99" or "-0.05%"), so NumberFormat is forced to Text ("@") before assignment to
# prevent Excel from auto-converting them into numeric/percentage values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value2 = "302.99"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value2 = "-0.05%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value2 = "37.19"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value2 = "6.25%"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value2 = "-3.85%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value2 = "0.07824"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value2 = "0.89%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value2 = "2.205"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value2 = "-3.67%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value2 = "8.016"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value2 = "-0.28%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value2 = "4.036"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value2 = "0.71%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value2 = "0.9142"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value2 = "-1.56%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value2 = "0.09758"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value2 = "-3.67%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value2 = "0.1887"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value2 = "3.08%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value2 = "0.08662"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value2 = "0.39%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value2 = "0.03569"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value2 = "3.33%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value2 = "0.09962"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value2 = "0.45%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value2 = "0.001486"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value2 = "0.11%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value2 = "0.005650"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value2 = "-2.54%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value2 = "3.462"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value2 = "-1.29%"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value2 = "7.15%"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value2 = "-2.19%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value2 = "4.782"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value2 = "3.91%"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value2 = "-2.19%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value2 = "0.04633"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value2 = "0.39%"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value2 = "8.33%"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value2 = "-7.88%"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value2 = "38.63%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value2 = "0.01774"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value2 = "0.38%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value2 = "0.04743"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value2 = "0.41%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value2 = "0.008063"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value2 = "5.77%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value2 = "0.1392"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value2 = "-0.98%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value2 = "0.007669"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value2 = "8.02%"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value2 = "-2.58%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value2 = "0.009870"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value2 = "7.11%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value2 = "0.00006103"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value2 = "2.90%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value2 = "-0.43%"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value2 = "-0.74%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value2 = "0.00002101"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value2 = "-0.43%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value2 = "0.0002001"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value2 = "-0.43%"
